$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "view" weights: the float slider expects raw fractional values
# (e.g. 0.05 for a 5% view) instead of the previously mis-scaled numbers,
# so the weights read in from the slider are corrected here.
$ws.Range("B2").Value = 0.2454406141024289
$ws.Range("C2").Value = 0.245440628078421

$ws.Range("B3").Value = 0.1815773618901335
$ws.Range("C3").Value = 0.1815773450933061

$ws.Range("B4").Value = 0.09883777827451906
$ws.Range("C4").Value = 0.09883777827451899

$ws.Range("B5").Value = 0.09883777827451906
$ws.Range("C5").Value = 0.09883777827451901

$ws.Range("B6").Value = 0.09883777827451906
$ws.Range("C6").Value = 0.09883777827451899

$ws.Range("B7").Value = 0.11995695347692
$ws.Range("C7").Value = 0.1199569551788409

$ws.Range("B8").Value = 0.1565117357069661
$ws.Range("C8").Value = 0.1565117368258783
